$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers (new date columns) - match header style of existing header cells
$ws.Range("AB1").Value = "31/03/2024"
$ws.Range("AC1").Value = "30/06/2024"
$ws.Range("AB1:AC1").Font.Bold = $true
$ws.Range("AB1:AC1").HorizontalAlignment = -4108
$ws.Range("AB1:AC1").VerticalAlignment = -4160
$ws.Range("AB1:AC1").Borders.LineStyle = 1

$ws.Range("AB2").Value = 8767729.664000001
$ws.Range("AC2").Value = 8583604.223999999
$ws.Range("AB3").Value = 4272453.888
$ws.Range("AC3").Value = 4191468.032
$ws.Range("AB4").Value = 232516
$ws.Range("AC4").Value = 108208
$ws.Range("AB5").Value = 0
$ws.Range("AC5").Value = 0
$ws.Range("AB6").Value = 882542.976
$ws.Range("AC6").Value = 1052702.016
$ws.Range("AB7").Value = 2886881.024
$ws.Range("AC7").Value = 2798281.984
$ws.Range("AB8").Value = 0
$ws.Range("AC8").Value = 0
$ws.Range("AB9").Value = 270513.984
$ws.Range("AC9").Value = 232276
$ws.Range("AB10").Value = 0
$ws.Range("AC10").Value = 0
$ws.Range("AB11").Value = 0
$ws.Range("AC11").Value = 0
$ws.Range("AB12").Value = 1388571.008
$ws.Range("AC12").Value = 1373102.976
$ws.Range("AB13").Value = 0
$ws.Range("AC13").Value = 0
$ws.Range("AB14").Value = 0
$ws.Range("AC14").Value = 0
$ws.Range("AB15").Value = 793
$ws.Range("AC15").Value = 809
$ws.Range("AB16").Value = 0
$ws.Range("AC16").Value = 0
$ws.Range("AB17").Value = 0
$ws.Range("AC17").Value = 0
$ws.Range("AB18").Value = 0
$ws.Range("AC18").Value = 0
$ws.Range("AB19").Value = 556846.976
$ws.Range("AC19").Value = 564350.0159999999
$ws.Range("AB20").Value = 0
$ws.Range("AC20").Value = 0
$ws.Range("AB21").Value = 0
$ws.Range("AC21").Value = 0
$ws.Range("AB22").Value = 78952
$ws.Range("AC22").Value = 80446
$ws.Range("AB23").Value = 920848
$ws.Range("AC23").Value = 914363.008
$ws.Range("AB24").Value = 2106904.96
$ws.Range("AC24").Value = 2024224
$ws.Range("AB25").Value = 0
$ws.Range("AC25").Value = 0
$ws.Range("AB26").Value = 8767729.664000001
$ws.Range("AC26").Value = 8583604.223999999
$ws.Range("AB27").Value = 3088506.88
$ws.Range("AC27").Value = 2811239.936
$ws.Range("AB28").Value = 194592
$ws.Range("AC28").Value = 224288.992
$ws.Range("AB29").Value = 1749330.944
$ws.Range("AC29").Value = 1775297.024
$ws.Range("AB30").Value = 118181
$ws.Range("AC30").Value = 52756
$ws.Range("AB31").Value = 415936
$ws.Range("AC31").Value = 138731.008
$ws.Range("AB32").Value = 0
$ws.Range("AC32").Value = 0
$ws.Range("AB33").Value = 0
$ws.Range("AC33").Value = 0
$ws.Range("AB34").Value = 610467.008
$ws.Range("AC34").Value = 620166.976
$ws.Range("AB35").Value = 0
$ws.Range("AC35").Value = 0
$ws.Range("AB36").Value = 0
$ws.Range("AC36").Value = 0
$ws.Range("AB37").Value = 3069853.952
$ws.Range("AC37").Value = 3127474.944
$ws.Range("AB38").Value = 1178370.944
$ws.Range("AC38").Value = 1330660.992
$ws.Range("AB39").Value = 0
$ws.Range("AC39").Value = 0
$ws.Range("AB40").Value = 1813730.944
$ws.Range("AC40").Value = 1731090.048
$ws.Range("AB41").Value = 0
$ws.Range("AC41").Value = 0
$ws.Range("AB42").Value = 0
$ws.Range("AC42").Value = 0
$ws.Range("AB43").Value = 77752
$ws.Range("AC43").Value = 65724
$ws.Range("AB44").Value = 0
$ws.Range("AC44").Value = 0
$ws.Range("AB45").Value = 0
$ws.Range("AC45").Value = 0
$ws.Range("AB46").Value = 7420
$ws.Range("AC46").Value = 7337
$ws.Range("AB47").Value = 2601949.088
$ws.Range("AC47").Value = 2637552.088
$ws.Range("AB48").Value = 1721858.048
$ws.Range("AC48").Value = 1721858.048
$ws.Range("AB49").Value = 379620
$ws.Range("AC49").Value = 383127.008
$ws.Range("AB50").Value = 0
$ws.Range("AC50").Value = 0
$ws.Range("AB51").Value = 537347.008
$ws.Range("AC51").Value = 536828.992
$ws.Range("AB52").Value = -36876
$ws.Range("AC52").Value = -4262
$ws.Range("AB53").Value = 0
$ws.Range("AC53").Value = 0
$ws.Range("AB54").Value = 0
$ws.Range("AC54").Value = 0
$ws.Range("AB55").Value = 0
$ws.Range("AC55").Value = 0
$ws.Range("AB56").Value = 0
$ws.Range("AC56").Value = 0
# Row 57: blank placeholder cells (AB57, AC57) - left empty intentionally
# Row 58: blank placeholder cells (AB58, AC58) - left empty intentionally
$ws.Range("AB59").Value = 2882500.096
$ws.Range("AC59").Value = 3138831.104
$ws.Range("AB60").Value = -1978946.944
$ws.Range("AC60").Value = -2115639.04
$ws.Range("AB61").Value = 903553.024
$ws.Range("AC61").Value = 1023192
$ws.Range("AB62").Value = -728540.992
$ws.Range("AC62").Value = -764489.024
$ws.Range("AB63").Value = -90555
$ws.Range("AC63").Value = -94231
$ws.Range("AB64").Value = 0
$ws.Range("AC64").Value = 0
$ws.Range("AB65").Value = 579
$ws.Range("AC65").Value = 652
$ws.Range("AB66").Value = -1030
$ws.Range("AC66").Value = -232
$ws.Range("AB67").Value = 1778
$ws.Range("AC67").Value = 1493
$ws.Range("AB68").Value = -157286
$ws.Range("AC68").Value = -141356.992
$ws.Range("AB69").Value = 13429
$ws.Range("AC69").Value = 58415
$ws.Range("AB70").Value = -170715.008
$ws.Range("AC70").Value = -199772
# Row 71: blank placeholder cells (AB71, AC71) - left empty intentionally
# Row 72: blank placeholder cells (AB72, AC72) - left empty intentionally
# Row 73: blank placeholder cells (AB73, AC73) - left empty intentionally
$ws.Range("AB74").Value = -71502
$ws.Range("AC74").Value = 25028
$ws.Range("AB75").Value = 0
$ws.Range("AC75").Value = 0
$ws.Range("AB76").Value = 34529
$ws.Range("AC76").Value = 7503
# Row 77: blank placeholder cells (AB77, AC77) - left empty intentionally
# Row 78: blank placeholder cells (AB78, AC78) - left empty intentionally
$ws.Range("AB79").Value = 97
$ws.Range("AC79").Value = 83
$ws.Range("AB80").Value = -36876
$ws.Range("AC80").Value = 32614

Write-Output "Applied PGMN3 Q1/Q2 2024 balance columns (AB:AC)"